$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.077.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.72%  '
$ws.Range("D3").Value = "'1.597.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.33%  '
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("D5").Value = "'1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("D6").Value = "'301.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.34%  '
$ws.Range("D7").Value = "'0.3771"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.35%  '
$ws.Range("D8").Value = "'0.3653"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.75%  '
$ws.Range("D9").Value = "'47.77"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.85%  '
$ws.Range("D10").Value = "'1.003"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("D11").Value = "'1.275"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.88%  '
$ws.Range("D12").Value = "'0.08070"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.38%  '
$ws.Range("D13").Value = "'22.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.20%  '
$ws.Range("D14").Value = "'6.623"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.98%  '
$ws.Range("D15").Value = "'7.629"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.82%  '
$ws.Range("D16").Value = "'0.00001265"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.88%  '
$ws.Range("D17").Value = "'1.594.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.48%  '
$ws.Range("D18").Value = "'91.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.21%  '
$ws.Range("D19").Value = "'0.06791"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.65%  '
$ws.Range("D20").Value = "'18.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.07%  '
$ws.Range("D21").Value = "'6.583"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.43%  '
$ws.Range("D22").Value = "'1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.36%  '
$ws.Range("D23").Value = "'13.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.44%  '
$ws.Range("D24").Value = "'23.085.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.70%  '
$ws.Range("D25").Value = "'2.363"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.56%  '
$ws.Range("D26").Value = "'2.884"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.38%  '
$ws.Range("D27").Value = "'21.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.41%  '
$ws.Range("D28").Value = "'151.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.06%  '
$ws.Range("D29").Value = "'5.252"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.28%  '
$ws.Range("D30").Value = "'132.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.19%  '
$ws.Range("D31").Value = "'2.437"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.83%  '
$ws.Range("D32").Value = "'7.116"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.17%  '
$ws.Range("D33").Value = "'1.772.06"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.45%  '
$ws.Range("D34").Value = "'0.9852"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.94%  '
$ws.Range("D35").Value = "'0.07708"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.21%  '
$ws.Range("D36").Value = "'0.02780"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.37%  '
$ws.Range("D37").Value = "'6.296"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.13%  '
$ws.Range("D38").Value = "'0.2541"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.28%  '
$ws.Range("D39").Value = "'0.08876"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.76%  '
$ws.Range("D40").Value = "'10.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.41%  '
$ws.Range("D41").Value = "'1.396"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.66%  '
$ws.Range("D42").Value = "'0.7145"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.55%  '
$ws.Range("D43").Value = "'12.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.29%  '
$ws.Range("D44").Value = "'15.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.15%  '
$ws.Range("D45").Value = "'0.6627"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.52%  '
$ws.Range("D46").Value = "'2.309"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.81%  '
$ws.Range("D47").Value = "'1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.22%  '
$ws.Range("D48").Value = "'3.963"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.65%  '
$ws.Range("D49").Value = "'131.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.89%  '
$ws.Range("D50").Value = "'0.07980"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.02%  '
$ws.Range("D51").Value = "'1.170"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.55%  '
